$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"30.39114566666667"
$ws.Range("H2").Value = [double]"91.17343700000001"
$ws.Range("I2").Value = [double]"0.2485034818803364"
$ws.Range("J2").Value = [double]"0.2485034818803363"
$ws.Range("K2").Value = [double]"1"
$ws.Range("L2").Value = [double]"0.3333333333333333"
$ws.Range("M2").Value = [double]"0.1145113333333333"
$ws.Range("N2").Value = [double]"0.343534"
$ws.Range("O2").Value = [double]"0.001785365609625045"
$ws.Range("P2").Value = [double]"0.001785365609625044"
$ws.Range("Q2").Value = [double]"3.480130611817556"
$ws.Range("R2").Value = [double]"31.321175506358"
$ws.Range("S2").Value = [double]"0.0004436695704212329"
$ws.Range("T2").Value = [double]"0.0004436695704212328"

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"30.39114566666667"
$ws.Range("H3").Value = [double]"91.17343700000001"
$ws.Range("I3").Value = [double]"0.2485034818803364"
$ws.Range("J3").Value = [double]"0.2485034818803363"
$ws.Range("K3").Value = [double]"3"
$ws.Range("L3").Value = [double]"1"
$ws.Range("M3").Value = [double]"0.467525"
$ws.Range("N3").Value = [double]"1.402575"
$ws.Range("O3").Value = [double]"0.007289261528465441"
$ws.Range("P3").Value = [double]"0.007289261528465441"
$ws.Range("Q3").Value = [double]"14.20862037780834"
$ws.Range("R3").Value = [double]"127.877583400275"
$ws.Range("S3").Value = [double]"0.001811406870160044"
$ws.Range("T3").Value = [double]"0.001811406870160044"

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"30.39114566666667"
$ws.Range("H4").Value = [double]"91.17343700000001"
$ws.Range("I4").Value = [double]"0.2485034818803364"
$ws.Range("J4").Value = [double]"0.2485034818803363"
$ws.Range("K4").Value = [double]"3"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"63.556834"
$ws.Range("N4").Value = [double]"190.670502"
$ws.Range("O4").Value = [double]"0.9909253728619096"
$ws.Range("P4").Value = [double]"0.9909253728619095"
$ws.Range("Q4").Value = [double]"1931.565000206153"
$ws.Range("R4").Value = [double]"17384.08500185538"
$ws.Range("S4").Value = [double]"0.2462484054397551"
$ws.Range("T4").Value = [double]"0.246248405439755"

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"45.91529066666667"
$ws.Range("H5").Value = [double]"137.745872"
$ws.Range("I5").Value = [double]"0.3754419042757282"
$ws.Range("J5").Value = [double]"0.3754419042757282"
$ws.Range("K5").Value = [double]"1"
$ws.Range("L5").Value = [double]"0.3333333333333333"
$ws.Range("M5").Value = [double]"0.1145113333333333"
$ws.Range("N5").Value = [double]"0.343534"
$ws.Range("O5").Value = [double]"0.001785365609625045"
$ws.Range("P5").Value = [double]"0.001785365609625044"
$ws.Range("Q5").Value = [double]"5.257821154627556"
$ws.Range("R5").Value = [double]"47.320390391648"
$ws.Range("S5").Value = [double]"0.0006703010643060231"
$ws.Range("T5").Value = [double]"0.000670301064306023"

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"45.91529066666667"
$ws.Range("H6").Value = [double]"137.745872"
$ws.Range("I6").Value = [double]"0.3754419042757282"
$ws.Range("J6").Value = [double]"0.3754419042757282"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"0.467525"
$ws.Range("N6").Value = [double]"1.402575"
$ws.Range("O6").Value = [double]"0.007289261528465441"
$ws.Range("P6").Value = [double]"0.007289261528465441"
$ws.Range("Q6").Value = [double]"21.46654626893334"
$ws.Range("R6").Value = [double]"193.1989164204001"
$ws.Range("S6").Value = [double]"0.00273669422901087"
$ws.Range("T6").Value = [double]"0.00273669422901087"

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"45.91529066666667"
$ws.Range("H7").Value = [double]"137.745872"
$ws.Range("I7").Value = [double]"0.3754419042757282"
$ws.Range("J7").Value = [double]"0.3754419042757282"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"63.556834"
$ws.Range("N7").Value = [double]"190.670502"
$ws.Range("O7").Value = [double]"0.9909253728619096"
$ws.Range("P7").Value = [double]"0.9909253728619095"
$ws.Range("Q7").Value = [double]"2918.230506963083"
$ws.Range("R7").Value = [double]"26264.07456266775"
$ws.Range("S7").Value = [double]"0.3720349089824113"
$ws.Range("T7").Value = [double]"0.3720349089824113"

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = [double]"3"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"2.332475"
$ws.Range("H8").Value = [double]"6.997425"
$ws.Range("I8").Value = [double]"0.01907227076123622"
$ws.Range("J8").Value = [double]"0.01907227076123622"
$ws.Range("K8").Value = [double]"1"
$ws.Range("L8").Value = [double]"0.3333333333333333"
$ws.Range("M8").Value = [double]"0.1145113333333333"
$ws.Range("N8").Value = [double]"0.343534"
$ws.Range("O8").Value = [double]"0.001785365609625045"
$ws.Range("P8").Value = [double]"0.001785365609625044"
$ws.Range("Q8").Value = [double]"0.2670948222166667"
$ws.Range("R8").Value = [double]"2.40385339995"
$ws.Range("S8").Value = [double]"3.405097631456842E-05"
$ws.Range("T8").Value = [double]"3.405097631456841E-05"

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = [double]"3"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"2.332475"
$ws.Range("H9").Value = [double]"6.997425"
$ws.Range("I9").Value = [double]"0.01907227076123622"
$ws.Range("J9").Value = [double]"0.01907227076123622"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"0.467525"
$ws.Range("N9").Value = [double]"1.402575"
$ws.Range("O9").Value = [double]"0.007289261528465441"
$ws.Range("P9").Value = [double]"0.007289261528465441"
$ws.Range("Q9").Value = [double]"1.090490374375"
$ws.Range("R9").Value = [double]"9.814413369375"
$ws.Range("S9").Value = [double]"0.0001390227695203555"
$ws.Range("T9").Value = [double]"0.0001390227695203555"

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = [double]"3"
$ws.Range("F10").Value = [double]"1"
$ws.Range("G10").Value = [double]"2.332475"
$ws.Range("H10").Value = [double]"6.997425"
$ws.Range("I10").Value = [double]"0.01907227076123622"
$ws.Range("J10").Value = [double]"0.01907227076123622"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"63.556834"
$ws.Range("N10").Value = [double]"190.670502"
$ws.Range("O10").Value = [double]"0.9909253728619096"
$ws.Range("P10").Value = [double]"0.9909253728619095"
$ws.Range("Q10").Value = [double]"148.24472638415"
$ws.Range("R10").Value = [double]"1334.20253745735"
$ws.Range("S10").Value = [double]"0.0188991970154013"
$ws.Range("T10").Value = [double]"0.01889919701540129"

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = [double]"3"
$ws.Range("F11").Value = [double]"1"
$ws.Range("G11").Value = [double]"43.657748"
$ws.Range("H11").Value = [double]"130.973244"
$ws.Range("I11").Value = [double]"0.3569823430826993"
$ws.Range("J11").Value = [double]"0.3569823430826993"
$ws.Range("K11").Value = [double]"1"
$ws.Range("L11").Value = [double]"0.3333333333333333"
$ws.Range("M11").Value = [double]"0.1145113333333333"
$ws.Range("N11").Value = [double]"0.343534"
$ws.Range("O11").Value = [double]"0.001785365609625045"
$ws.Range("P11").Value = [double]"0.001785365609625044"
$ws.Range("Q11").Value = [double]"4.999306933810667"
$ws.Range("R11").Value = [double]"44.993762404296"
$ws.Range("S11").Value = [double]"0.0006373439985832202"
$ws.Range("T11").Value = [double]"0.0006373439985832201"

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = [double]"3"
$ws.Range("F12").Value = [double]"1"
$ws.Range("G12").Value = [double]"43.657748"
$ws.Range("H12").Value = [double]"130.973244"
$ws.Range("I12").Value = [double]"0.3569823430826993"
$ws.Range("J12").Value = [double]"0.3569823430826993"
$ws.Range("K12").Value = [double]"3"
$ws.Range("L12").Value = [double]"1"
$ws.Range("M12").Value = [double]"0.467525"
$ws.Range("N12").Value = [double]"1.402575"
$ws.Range("O12").Value = [double]"0.007289261528465441"
$ws.Range("P12").Value = [double]"0.007289261528465441"
$ws.Range("Q12").Value = [double]"20.4110886337"
$ws.Range("R12").Value = [double]"183.6997977033"
$ws.Range("S12").Value = [double]"0.002602137659774171"
$ws.Range("T12").Value = [double]"0.002602137659774171"

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = [double]"3"
$ws.Range("F13").Value = [double]"1"
$ws.Range("G13").Value = [double]"43.657748"
$ws.Range("H13").Value = [double]"130.973244"
$ws.Range("I13").Value = [double]"0.3569823430826993"
$ws.Range("J13").Value = [double]"0.3569823430826993"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"63.556834"
$ws.Range("N13").Value = [double]"190.670502"
$ws.Range("O13").Value = [double]"0.9909253728619096"
$ws.Range("P13").Value = [double]"0.9909253728619095"
$ws.Range("Q13").Value = [double]"2774.748242449832"
$ws.Range("R13").Value = [double]"24972.73418204849"
$ws.Range("S13").Value = [double]"0.3537428614243419"
$ws.Range("T13").Value = [double]"0.3537428614243419"

# Remove old rows 14-17 (Resolving-Mac sending-cluster block superseded)
$ws.Rows("14:17").Delete()
